$d = $word.ActiveDocument

# The "ncombinations" paragraph currently holds ": <number>" in a single
# run/paragraph. Split it into two paragraphs: the first keeps ": " and a
# new paragraph (cloning the original paragraph's spacing) holds just the
# number, by replacing the space after the colon with a paragraph break.
$old = ": 134779202705884333073115306260788792805547398513295360"
$new = ": ^p134779202705884333073115306260788792805547398513295360"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target text to split into two paragraphs"
}
